$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-02-27 Thursday" "2025-02-28 Friday"

Replace-Text "296÷5=59, 1" "932÷3=310, 2"
Replace-Text "868÷2=434, 0" "233÷9=25, 8"
Replace-Text "501÷6=83, 3" "136÷3=45, 1"
Replace-Text "361÷7=51, 4" "334÷7=47, 5"
Replace-Text "159÷9=17, 6" "478÷5=95, 3"

Replace-Text "585÷9=65, 0" "549÷9=61, 0"
Replace-Text "338÷2=169, 0" "374÷6=62, 2"
Replace-Text "815÷7=116, 3" "679÷5=135, 4"
Replace-Text "590÷7=84, 2" "274÷4=68, 2"
Replace-Text "393÷7=56, 1" "693÷7=99, 0"

Replace-Text "992÷5=198, 2" "653÷2=326, 1"
Replace-Text "750÷6=125, 0" "264÷8=33, 0"
Replace-Text "890÷8=111, 2" "310÷2=155, 0"
Replace-Text "106÷7=15, 1" "945÷5=189, 0"
Replace-Text "870÷4=217, 2" "513÷9=57, 0"

Replace-Text "940÷6=156, 4" "425÷9=47, 2"
Replace-Text "866÷7=123, 5" "395÷4=98, 3"
Replace-Text "653÷6=108, 5" "145÷3=48, 1"
Replace-Text "477÷9=53, 0" "927÷6=154, 3"
Replace-Text "883÷3=294, 1" "544÷4=136, 0"

Replace-Text "735÷6=122, 3" "466÷4=116, 2"
Replace-Text "552÷9=61, 3" "558÷4=139, 2"
Replace-Text "914÷7=130, 4" "141÷8=17, 5"
Replace-Text "577÷4=144, 1" "830÷6=138, 2"
Replace-Text "692÷8=86, 4" "590÷7=84, 2"
